$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 15)
$ws.Range("A15").Value = 6025
$ws.Range("B15").Value = "FIKRI MUHAMMAD RIZAL"
$ws.Range("C15").Value = 18

# Update the active selection to match the new state (A15)
$ws.Range("A15").Select()
